$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '310.30'
Set-TextValue 'E2' '1.47%'
Set-TextValue 'D3' '35.61'
Set-TextValue 'E3' '-1.96%'
Set-TextValue 'D4' '5.103'
Set-TextValue 'E4' '0.81%'
Set-TextValue 'D5' '0.08215'
Set-TextValue 'E5' '4.08%'
Set-TextValue 'D6' '2.059'
Set-TextValue 'E6' '-3.34%'
Set-TextValue 'B7' 'KuCoinToken'
Set-TextValue 'C7' 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue 'D7' '7.933'
Set-TextValue 'E7' '-0.67%'
Set-TextValue 'B8' 'BTSEToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue 'D8' '2.986'
Set-TextValue 'E8' '12.24%'
Set-TextValue 'B9' 'MXToken'
Set-TextValue 'C9' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D9' '0.9264'
Set-TextValue 'E9' '-0.20%'
Set-TextValue 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D10' '0.1125'
Set-TextValue 'E10' '15.43%'
Set-TextValue 'B11' 'WazirX'
Set-TextValue 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.1920'
Set-TextValue 'E11' '2.89%'
Set-TextValue 'B12' 'MandalaExchangeToken'
Set-TextValue 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D12' '0.09281'
Set-TextValue 'E12' '3.18%'
Set-TextValue 'B13' 'BitrueCoin'
Set-TextValue 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.03577'
Set-TextValue 'E13' '-5.08%'
Set-TextValue 'B14' 'BitMartToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09907'
Set-TextValue 'E14' '0.14%'
Set-TextValue 'B15' 'BitForexToken'
Set-TextValue 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D15' '0.001439'
Set-TextValue 'E15' '0.13%'
Set-TextValue 'B16' 'TigerCash'
Set-TextValue 'C16' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D16' '0.005738'
Set-TextValue 'E16' '1.52%'
Set-TextValue 'B17' 'LEO'
Set-TextValue 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D17' '3.468'
Set-TextValue 'E17' '-0.02%'
Set-TextValue 'B18' 'GateToken'
Set-TextValue 'C18' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D18' '4.126'
Set-TextValue 'E18' '-0.74%'
Set-TextValue 'D19' '0.3428'
Set-TextValue 'E19' '0.17%'
Set-TextValue 'E20' '-0.28%'
Set-TextValue 'D21' '5.097'
Set-TextValue 'E21' '-0.68%'
Set-TextValue 'D23' '0.04539'
Set-TextValue 'E23' '-1.21%'
Set-TextValue 'E24' '-1.10%'
Set-TextValue 'D25' '0.004813'
Set-TextValue 'E25' '0.41%'
Set-TextValue 'D27' '0.0004446'
Set-TextValue 'E27' '-6.16%'
Set-TextValue 'D39' '0.01989'
Set-TextValue 'E39' '1.60%'
Set-TextValue 'D40' '0.04933'
Set-TextValue 'E40' '-0.24%'
Set-TextValue 'D41' '0.007678'
Set-TextValue 'E41' '-1.92%'
Set-TextValue 'D42' '0.01000'
Set-TextValue 'E42' '27.76%'
Set-TextValue 'D43' '0.1385'
Set-TextValue 'E43' '-0.78%'
Set-TextValue 'D44' '0.002122'
Set-TextValue 'E44' '-0.53%'
Set-TextValue 'E45' '2.94%'
Set-TextValue 'D46' '0.00006555'
Set-TextValue 'E46' '4.41%'
Set-TextValue 'E47' '-0.35%'
Set-TextValue 'D48' '61.63'
Set-TextValue 'E48' '19.23%'
Set-TextValue 'E49' '-21.40%'
Set-TextValue 'D50' '0.00002099'
Set-TextValue 'E50' '-0.35%'
Set-TextValue 'D51' '0.0001999'
Set-TextValue 'E51' '-0.35%'
